$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "JOSE GUEVARA ANAYA" row (row 16) -----------------------
# Deleting row 16 shifts the two MARIA SILVIA rows up:
#   old row17 (MARIA SILVIA / 2504, normal style)      -> new row16
#   old row18 (MARIA SILVIA / 2506, bottom-border style)-> new row17
# This naturally reproduces the correct per-row borders/styles without any
# manual style surgery.
$ws.Rows.Item(16).Delete()

# --- Turn the (now) row 17 back into the JOSE GUEVARA ANAYA record ----------
# (same values as the old row16, but now placed last, with the bottom-border
# row style that came along for free from the old row18).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "7384405"
$ws.Range("D17").Value = "JOSE GUEVARA ANAYA"
$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

# --- Update the summary figures ---------------------------------------------
# Valor Mora total: 168960 -> 110480 (one 58480 period removed)
$ws.Range("E11").Value = 110480
# Cant. Periodos: 2 -> 1 (only period 2504 remains)
$ws.Range("F13").Value = 1
